$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 3.767119
$ws.Cells.Item(2, 8).Value = 7.534238
$ws.Cells.Item(2, 9).Value = 0.1064531723628092
$ws.Cells.Item(2, 10).Value = 0.08010427729187669
$ws.Cells.Item(2, 13).Value = 1.273363333333333
$ws.Cells.Item(2, 14).Value = 3.82009
$ws.Cells.Item(2, 15).Value = 0.1769428433887536
$ws.Cells.Item(2, 16).Value = 0.1769428433887536
$ws.Cells.Item(2, 17).Value = 4.796911206903334
$ws.Cells.Item(2, 18).Value = 28.78146724142
$ws.Cells.Item(2, 19).Value = 0.01883612700562854
$ws.Cells.Item(2, 20).Value = 0.01417387859162583
# Row 3
$ws.Cells.Item(3, 7).Value = 3.767119
$ws.Cells.Item(3, 8).Value = 7.534238
$ws.Cells.Item(3, 9).Value = 0.1064531723628092
$ws.Cells.Item(3, 10).Value = 0.08010427729187669
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.815617
$ws.Cells.Item(3, 14).Value = 2.446851
$ws.Cells.Item(3, 15).Value = 0.1133357521128076
$ws.Cells.Item(3, 16).Value = 0.1133357521128076
$ws.Cells.Item(3, 17).Value = 3.072526297423
$ws.Cells.Item(3, 18).Value = 18.435157784538
$ws.Cells.Item(3, 19).Value = 0.01206495035453333
$ws.Cells.Item(3, 20).Value = 0.009078678514327736
# Row 4
$ws.Cells.Item(4, 7).Value = 3.767119
$ws.Cells.Item(4, 8).Value = 7.534238
$ws.Cells.Item(4, 9).Value = 0.1064531723628092
$ws.Cells.Item(4, 10).Value = 0.08010427729187669
$ws.Cells.Item(4, 13).Value = 5.107486666666667
$ws.Cells.Item(4, 14).Value = 15.32246
$ws.Cells.Item(4, 15).Value = 0.7097214044984388
$ws.Cells.Item(4, 16).Value = 0.7097214044984388
$ws.Cells.Item(4, 17).Value = 19.24051006424667
$ws.Cells.Item(4, 18).Value = 115.44306038548
$ws.Cells.Item(4, 19).Value = 0.07555209500264737
$ws.Cells.Item(4, 20).Value = 0.05685172018592312
# Row 5
$ws.Cells.Item(5, 9).Value = 0.01146655947899596
$ws.Cells.Item(5, 10).Value = 0.01294259869906215
$ws.Cells.Item(5, 13).Value = 1.273363333333333
$ws.Cells.Item(5, 14).Value = 3.82009
$ws.Cells.Item(5, 15).Value = 0.1769428433887536
$ws.Cells.Item(5, 16).Value = 0.1769428433887536
$ws.Cells.Item(5, 17).Value = 0.5166973087655556
$ws.Cells.Item(5, 18).Value = 4.650275778890001
$ws.Cells.Item(5, 19).Value = 0.00202892563809981
$ws.Cells.Item(5, 20).Value = 0.00229010021465164
# Row 6
$ws.Cells.Item(6, 9).Value = 0.01146655947899596
$ws.Cells.Item(6, 10).Value = 0.01294259869906215
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.815617
$ws.Cells.Item(6, 14).Value = 2.446851
$ws.Cells.Item(6, 15).Value = 0.1133357521128076
$ws.Cells.Item(6, 16).Value = 0.1133357521128076
$ws.Cells.Item(6, 17).Value = 0.3309559006856667
$ws.Cells.Item(6, 18).Value = 2.978603106171001
$ws.Cells.Item(6, 19).Value = 0.00129957114269825
$ws.Cells.Item(6, 20).Value = 0.001466859157852454
# Row 7
$ws.Cells.Item(7, 9).Value = 0.01146655947899596
$ws.Cells.Item(7, 10).Value = 0.01294259869906215
$ws.Cells.Item(7, 13).Value = 5.107486666666667
$ws.Cells.Item(7, 14).Value = 15.32246
$ws.Cells.Item(7, 15).Value = 0.7097214044984388
$ws.Cells.Item(7, 16).Value = 0.7097214044984388
$ws.Cells.Item(7, 17).Value = 2.072483592184445
$ws.Cells.Item(7, 18).Value = 18.65235232966
$ws.Cells.Item(7, 19).Value = 0.008138062698197897
$ws.Cells.Item(7, 20).Value = 0.009185639326558057
# Row 8
$ws.Cells.Item(8, 7).Value = 6.271924666666666
$ws.Cells.Item(8, 8).Value = 18.815774
$ws.Cells.Item(8, 9).Value = 0.1772352499581833
$ws.Cells.Item(8, 10).Value = 0.2000499556766435
$ws.Cells.Item(8, 13).Value = 1.273363333333333
$ws.Cells.Item(8, 14).Value = 3.82009
$ws.Cells.Item(8, 15).Value = 0.1769428433887536
$ws.Cells.Item(8, 16).Value = 0.1769428433887536
$ws.Cells.Item(8, 17).Value = 7.986438899962222
$ws.Cells.Item(8, 18).Value = 71.87795009966
$ws.Cells.Item(8, 19).Value = 0.03136050907631742
$ws.Cells.Item(8, 20).Value = 0.03539740797721944
# Row 9
$ws.Cells.Item(9, 7).Value = 6.271924666666666
$ws.Cells.Item(9, 8).Value = 18.815774
$ws.Cells.Item(9, 9).Value = 0.1772352499581833
$ws.Cells.Item(9, 10).Value = 0.2000499556766435
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.815617
$ws.Cells.Item(9, 14).Value = 2.446851
$ws.Cells.Item(9, 15).Value = 0.1133357521128076
$ws.Cells.Item(9, 16).Value = 0.1133357521128076
$ws.Cells.Item(9, 17).Value = 5.115488380852667
$ws.Cells.Item(9, 18).Value = 46.03939542767399
$ws.Cells.Item(9, 19).Value = 0.02008709035491215
$ws.Cells.Item(9, 20).Value = 0.02267281218674621
# Row 10
$ws.Cells.Item(10, 7).Value = 6.271924666666666
$ws.Cells.Item(10, 8).Value = 18.815774
$ws.Cells.Item(10, 9).Value = 0.1772352499581833
$ws.Cells.Item(10, 10).Value = 0.2000499556766435
$ws.Cells.Item(10, 13).Value = 5.107486666666667
$ws.Cells.Item(10, 14).Value = 15.32246
$ws.Cells.Item(10, 15).Value = 0.7097214044984388
$ws.Cells.Item(10, 16).Value = 0.7097214044984388
$ws.Cells.Item(10, 17).Value = 32.03377160933778
$ws.Cells.Item(10, 18).Value = 288.30394448404
$ws.Cells.Item(10, 19).Value = 0.1257876505269537
$ws.Cells.Item(10, 20).Value = 0.1419797355126779
# Row 11
$ws.Cells.Item(11, 7).Value = 8.340211500000001
$ws.Cells.Item(11, 8).Value = 16.680423
$ws.Cells.Item(11, 9).Value = 0.2356819554550265
$ws.Cells.Item(11, 10).Value = 0.1773468304741365
$ws.Cells.Item(11, 13).Value = 1.273363333333333
$ws.Cells.Item(11, 14).Value = 3.82009
$ws.Cells.Item(11, 15).Value = 0.1769428433887536
$ws.Cells.Item(11, 16).Value = 0.1769428433887536
$ws.Cells.Item(11, 17).Value = 10.620119516345
$ws.Cells.Item(11, 18).Value = 63.72071709807001
$ws.Cells.Item(11, 19).Value = 0.04170223533363394
$ws.Cells.Item(11, 20).Value = 0.03138025245007697
# Row 12
$ws.Cells.Item(12, 7).Value = 8.340211500000001
$ws.Cells.Item(12, 8).Value = 16.680423
$ws.Cells.Item(12, 9).Value = 0.2356819554550265
$ws.Cells.Item(12, 10).Value = 0.1773468304741365
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.815617
$ws.Cells.Item(12, 14).Value = 2.446851
$ws.Cells.Item(12, 15).Value = 0.1133357521128076
$ws.Cells.Item(12, 16).Value = 0.1133357521128076
$ws.Cells.Item(12, 17).Value = 6.802418282995501
$ws.Cells.Item(12, 18).Value = 40.814509697973
$ws.Cells.Item(12, 19).Value = 0.02671119168091263
$ws.Cells.Item(12, 20).Value = 0.02009973641660885
# Row 13
$ws.Cells.Item(13, 7).Value = 8.340211500000001
$ws.Cells.Item(13, 8).Value = 16.680423
$ws.Cells.Item(13, 9).Value = 0.2356819554550265
$ws.Cells.Item(13, 10).Value = 0.1773468304741365
$ws.Cells.Item(13, 13).Value = 5.107486666666667
$ws.Cells.Item(13, 14).Value = 15.32246
$ws.Cells.Item(13, 15).Value = 0.7097214044984388
$ws.Cells.Item(13, 16).Value = 0.7097214044984388
$ws.Cells.Item(13, 17).Value = 42.59751903343001
$ws.Cells.Item(13, 18).Value = 255.5851142005801
$ws.Cells.Item(13, 19).Value = 0.1672685284404799
$ws.Cells.Item(13, 20).Value = 0.1258668416074507
# Row 14
$ws.Cells.Item(14, 7).Value = 4.504435666666667
$ws.Cells.Item(14, 8).Value = 13.513307
$ws.Cells.Item(14, 9).Value = 0.1272886432366092
$ws.Cells.Item(14, 10).Value = 0.1436739443402582
$ws.Cells.Item(14, 13).Value = 1.273363333333333
$ws.Cells.Item(14, 14).Value = 3.82009
$ws.Cells.Item(14, 15).Value = 0.1769428433887536
$ws.Cells.Item(14, 16).Value = 0.1769428433887536
$ws.Cells.Item(14, 17).Value = 5.735783215292223
$ws.Cells.Item(14, 18).Value = 51.62204893763001
$ws.Cells.Item(14, 19).Value = 0.02252281446538228
$ws.Cells.Item(14, 20).Value = 0.02542207623244281
# Row 15
$ws.Cells.Item(15, 7).Value = 4.504435666666667
$ws.Cells.Item(15, 8).Value = 13.513307
$ws.Cells.Item(15, 9).Value = 0.1272886432366092
$ws.Cells.Item(15, 10).Value = 0.1436739443402582
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.815617
$ws.Cells.Item(15, 14).Value = 2.446851
$ws.Cells.Item(15, 15).Value = 0.1133357521128076
$ws.Cells.Item(15, 16).Value = 0.1133357521128076
$ws.Cells.Item(15, 17).Value = 3.673894305139667
$ws.Cells.Item(15, 18).Value = 33.06504874625701
$ws.Cells.Item(15, 19).Value = 0.01442635411663994
$ws.Cells.Item(15, 20).Value = 0.01628339454081681
# Row 16
$ws.Cells.Item(16, 7).Value = 4.504435666666667
$ws.Cells.Item(16, 8).Value = 13.513307
$ws.Cells.Item(16, 9).Value = 0.1272886432366092
$ws.Cells.Item(16, 10).Value = 0.1436739443402582
$ws.Cells.Item(16, 13).Value = 5.107486666666667
$ws.Cells.Item(16, 14).Value = 15.32246
$ws.Cells.Item(16, 15).Value = 0.7097214044984388
$ws.Cells.Item(16, 16).Value = 0.7097214044984388
$ws.Cells.Item(16, 17).Value = 23.00634510835778
$ws.Cells.Item(16, 18).Value = 207.05710597522
$ws.Cells.Item(16, 19).Value = 0.09033947465458701
$ws.Cells.Item(16, 20).Value = 0.1019684735669986
# Row 17
$ws.Cells.Item(17, 7).Value = 12.09810466666667
$ws.Cells.Item(17, 8).Value = 36.294314
$ws.Cells.Item(17, 9).Value = 0.3418744195083758
$ws.Cells.Item(17, 10).Value = 0.3858823935180229
$ws.Cells.Item(17, 13).Value = 1.273363333333333
$ws.Cells.Item(17, 14).Value = 3.82009
$ws.Cells.Item(17, 15).Value = 0.1769428433887536
$ws.Cells.Item(17, 16).Value = 0.1769428433887536
$ws.Cells.Item(17, 17).Value = 15.40528288536222
$ws.Cells.Item(17, 18).Value = 138.64754596826
$ws.Cells.Item(17, 19).Value = 0.06049223186969158
$ws.Cells.Item(17, 20).Value = 0.06827912792273691
# Row 18
$ws.Cells.Item(18, 7).Value = 12.09810466666667
$ws.Cells.Item(18, 8).Value = 36.294314
$ws.Cells.Item(18, 9).Value = 0.3418744195083758
$ws.Cells.Item(18, 10).Value = 0.3858823935180229
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.815617
$ws.Cells.Item(18, 14).Value = 2.446851
$ws.Cells.Item(18, 15).Value = 0.1133357521128076
$ws.Cells.Item(18, 16).Value = 0.1133357521128076
$ws.Cells.Item(18, 17).Value = 9.867419833912667
$ws.Cells.Item(18, 18).Value = 88.806778505214
$ws.Cells.Item(18, 19).Value = 0.03874659446311127
$ws.Cells.Item(18, 20).Value = 0.04373427129645551
# Row 19
$ws.Cells.Item(19, 7).Value = 12.09810466666667
$ws.Cells.Item(19, 8).Value = 36.294314
$ws.Cells.Item(19, 9).Value = 0.3418744195083758
$ws.Cells.Item(19, 10).Value = 0.3858823935180229
$ws.Cells.Item(19, 13).Value = 5.107486666666667
$ws.Cells.Item(19, 14).Value = 15.32246
$ws.Cells.Item(19, 15).Value = 0.7097214044984388
$ws.Cells.Item(19, 16).Value = 0.7097214044984388
$ws.Cells.Item(19, 17).Value = 61.79090827693778
$ws.Cells.Item(19, 18).Value = 556.1181744924401
$ws.Cells.Item(19, 19).Value = 0.242635593175573
$ws.Cells.Item(19, 20).Value = 0.2738689942988305
